$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37.84278171914716
$ws.Range("C2").Value = 38.74047329475789
$ws.Range("D2").Value = 36.93294442240585
$ws.Range("E2").Value = 37.84294501151045
$ws.Range("F2").Value = 37.94553809617637
$ws.Range("G2").Value = 37.18388087321012
$ws.Range("J2").Value = 38.5495336601405
$ws.Range("K2").Value = 37.90966494870826
$ws.Range("L2").Value = 37.84710463969296
$ws.Range("M2").Value = 37.05173872635863
$ws.Range("N2").Value = 19.8418013726711
$ws.Range("O2").Value = 32.53665438233396
$ws.Range("P2").Value = 41.39802372017179
